$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("N2").Value = "OPSET"
$ws.Range("N5").Value = "funct6"

# Add new rows of funct6 values under column E (store as text, preserving leading zeros)
$ws.Range("E9:E11").NumberFormat = "@"
$ws.Range("E9").Value = "000011"
$ws.Range("E10").Value = "000100"
$ws.Range("E11").Value = "000101"

# Reset the number format back to the sheet's default style (General) without
# creating a new custom number-format style entry, by pasting formats from a
# default-styled cell.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E9:E11").PasteSpecial(-4122) | Out-Null
